$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "-0.08"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "-0.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.15"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "0.15"
$ws.Range("F2").Style = "Normal"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "0.06"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "-0.06"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.08"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.16"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "-0.19*"
$ws.Range("G4").Style = "Normal"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "0.14"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "-0.07"
$ws.Range("C5").Style = "Normal"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "-0.23**"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "0.02"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "-0.01"
$ws.Range("G6").Style = "Normal"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "0.09"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "-0.26**"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.02"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.32***"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "-0.04"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "0.08"
$ws.Range("G7").Style = "Normal"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "-0.15"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "-0.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.17"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "0.09"
$ws.Range("B9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.06"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "-0.19"
$ws.Range("F9").Style = "Normal"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "-0.05"
$ws.Range("B10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "-0.1"
$ws.Range("F10").Style = "Normal"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "-0.11"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "-0.12"
$ws.Range("C11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.14"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "0.14"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "-0.01"
$ws.Range("G11").Style = "Normal"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "-0.15"
$ws.Range("C12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.08"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.02"
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "-0.05"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "0.07"
$ws.Range("C13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.05"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "0.04"
$ws.Range("G13").Style = "Normal"
